$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: price (column D) and 1h volume change (column E) per coin row.
# Plain-looking-number price strings are written via a Text-format + paste-formats-back
# round trip so they stay text (matching the original inlineStr cells) instead of
# being auto-coerced to numbers by Excel.
$textUpdates = @(
    @{Ref="D2"; Val="37.400.43"},
    @{Ref="D3"; Val="2.048.91"},
    @{Ref="E3"; Val="  -1.40%  "},
    @{Ref="E4"; Val="  -0.01%  "},
    @{Ref="E5"; Val="  -1.63%  "},
    @{Ref="E6"; Val="  -1.37%  "},
    @{Ref="E7"; Val="  +0.06%  "},
    @{Ref="E8"; Val="  -2.27%  "},
    @{Ref="E9"; Val="  -1.41%  "},
    @{Ref="E10"; Val="  +0.56%  "},
    @{Ref="E11"; Val="  -1.78%  "},
    @{Ref="E12"; Val="  -0.55%  "},
    @{Ref="D13"; Val="2.351.21"},
    @{Ref="E13"; Val="  -1.32%  "},
    @{Ref="E14"; Val="  -1.74%  "},
    @{Ref="E15"; Val="  -2.85%  "},
    @{Ref="E16"; Val="  -0.28%  "},
    @{Ref="D17"; Val="2.034.83"},
    @{Ref="E17"; Val="  -1.88%  "},
    @{Ref="D18"; Val="37.281.76"},
    @{Ref="E18"; Val="  -1.20%  "},
    @{Ref="E19"; Val="  -0.85%  "},
    @{Ref="E20"; Val="  -2.58%  "},
    @{Ref="D21"; Val="0.0₃0828"},
    @{Ref="E21"; Val="  -1.31%  "},
    @{Ref="E22"; Val="  -1.35%  "},
    @{Ref="E23"; Val="  +0.09%  "},
    @{Ref="E24"; Val="  +0.18%  "},
    @{Ref="E25"; Val="  -3.82%  "},
    @{Ref="E26"; Val="  -0.69%  "},
    @{Ref="E27"; Val="  -2.14%  "},
    @{Ref="E28"; Val="  -5.85%  "},
    @{Ref="E29"; Val="  -2.48%  "},
    @{Ref="E30"; Val="  -3.07%  "},
    @{Ref="E31"; Val="  -2.22%  "},
    @{Ref="E32"; Val="  -3.62%  "},
    @{Ref="E33"; Val="  -2.59%  "},
    @{Ref="E34"; Val="  -2.11%  "},
    @{Ref="E35"; Val="  -1.11%  "},
    @{Ref="E36"; Val="  +1.36%  "},
    @{Ref="E37"; Val="  +0.13%  "},
    @{Ref="E38"; Val="  -4.41%  "},
    @{Ref="E39"; Val="  -1.65%  "},
    @{Ref="E40"; Val="  -4.94%  "},
    @{Ref="E41"; Val="  +2.22%  "},
    @{Ref="E42"; Val="  -1.24%  "},
    @{Ref="D43"; Val="1.473.70"},
    @{Ref="E43"; Val="  +1.77%  "},
    @{Ref="E44"; Val="  -3.25%  "},
    @{Ref="E45"; Val="  -4.81%  "},
    @{Ref="E46"; Val="  +0.97%  "},
    @{Ref="E47"; Val="  -3.89%  "},
    @{Ref="E48"; Val="  -4.80%  "},
    @{Ref="E49"; Val="  -2.39%  "},
    @{Ref="E50"; Val="  -2.02%  "},
    @{Ref="D51"; Val="2.240.38"},
    @{Ref="E51"; Val="  -1.28%  "}
)

foreach ($u in $textUpdates) {
    $ws.Range($u.Ref).Value = $u.Val
}

# Price cells whose new text looks like a plain number (e.g. "1.00", "5.32") need to be
# forced to Text so Excel does not reinterpret them as numeric values.
$plainNumberPrices = @(
    @{Ref="D4"; Val="1.00"},
    @{Ref="D5"; Val="229.23"},
    @{Ref="D8"; Val="57.16"},
    @{Ref="D9"; Val="0.387"},
    @{Ref="D10"; Val="0.0790"},
    @{Ref="D12"; Val="14.78"},
    @{Ref="D14"; Val="20.71"},
    @{Ref="D15"; Val="0.759"},
    @{Ref="D16"; Val="5.33"},
    @{Ref="D19"; Val="6.09"},
    @{Ref="D22"; Val="226.22"},
    @{Ref="D26"; Val="9.67"},
    @{Ref="D27"; Val="168.33"},
    @{Ref="D29"; Val="18.98"},
    @{Ref="D34"; Val="4.57"},
    @{Ref="D35"; Val="2.43"},
    @{Ref="D37"; Val="1.00"},
    @{Ref="D39"; Val="5.32"},
    @{Ref="D41"; Val="17.26"},
    @{Ref="D44"; Val="0.0942"},
    @{Ref="D45"; Val="96.33"},
    @{Ref="D48"; Val="3.92"},
    @{Ref="D49"; Val="7.16"}
)

foreach ($u in $plainNumberPrices) {
    $cell = $ws.Range($u.Ref)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Val
    # Restore the original (default/no special number format) cell style by pasting
    # formats only from an untouched donor cell, so only the value changes.
    $ws.Range("D6").Copy()
    $cell.PasteSpecial(-4122)
}
$excel.CutCopyMode = $false
